$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 569 (shifts rows 569..end down by one)
$ws.Rows.Item(569).Insert()

# New row 569 becomes a duplicate of the (still unshifted) row 568's original data
$ws.Cells.Item(569, 1).Value = 3
$ws.Cells.Item(569, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(569, 3).Value = "Coquimbo"
$ws.Cells.Item(569, 4).Value = 45005
$ws.Cells.Item(569, 5).Value = 5
$ws.Cells.Item(569, 6).Value = 100112037
$ws.Cells.Item(569, 7).Value = "Cebollín"
$ws.Cells.Item(569, 8).Value = "Sin especificar"
$ws.Cells.Item(569, 9).Value = "Primera"
$ws.Cells.Item(569, 10).Value = 250
$ws.Cells.Item(569, 11).Value = 3500
$ws.Cells.Item(569, 12).Value = 4000
$ws.Cells.Item(569, 13).Value = 3740
$ws.Cells.Item(569, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(569, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(569, 16).Value = 104
$ws.Cells.Item(569, 17).Value = 36
$ws.Cells.Item(569, 18).Value = "Hortaliza"

# Copy the date number format from row 568's date cell to the new row 569 date cell
$ws.Cells.Item(568, 4).Copy()
$ws.Cells.Item(569, 4).PasteSpecial(-4122)
$ws.Cells.Item(569, 4).Value = 45005

# Update row 568's Fecha (D) and Precio promedio ponderado (M) with the new values
$ws.Cells.Item(568, 4).Value = 45015
$ws.Cells.Item(568, 13).Value = 3760
